$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Collapse the three CORE COMPETENCIES paragraphs into a single,
#    much shorter summary paragraph.
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "Product Management & Strategy • Technical Product Development • Platform & Infrastructure"

$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)
$delRange = $d.Range($p7.Range.Start, $p8.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------
# 2. Append a new "TECHNICAL SKILLS" section at the very end of the
#    document, containing the detailed competency breakdown that used
#    to live in CORE COMPETENCIES.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$headingPara.Style = "Heading2"
$headingPara.Range.Text = "TECHNICAL SKILLS"

$headingPara.Range.InsertParagraphAfter()
$line1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$line1.Style = "Normal"
$line1.Range.Text = "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development; Product Strategy; Stakeholder Management; Product Analytics"

$line1.Range.InsertParagraphAfter()
$line2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$line2.Style = "Normal"
$line2.Range.Text = "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; Database Design; API Development; DevOps & Deployment; System Integration"

$line2.Range.InsertParagraphAfter()
$line3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$line3.Style = "Normal"
$line3.Range.Text = "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Real-time Systems; Security & Compliance; Monitoring & Analytics; Documentation & Training"
